# Burn Report - "Syncing from share - 4/17" update
# Adds five new burn-log entries (rows 60-64) dated 4/17/2018 (serial 43207)
# for: GUI Movement Map (Cody), encoder test bench (Peter), LED integration
# (Brandon), Led VHDL (Brandon), and C function verification (Zack).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 60: GUI - Movement Map -------------------------------------------
$ws.Range("A60").Value = 43207
$ws.Range("B60").Value = "GUI - Movement Map"
$ws.Range("C60").Value = 2
$ws.Range("G60").Value = 1

# --- Row 61: Finished encoder test bench. ----------------------------------
# Hours (C61) is entered as text "1" in the source workbook, so force a
# text-formatted round-trip instead of letting it coerce to a number.
$ws.Range("A61").Value = 43207
$ws.Range("B61").Value = "Finished encoder test bench."
$ws.Range("C61").NumberFormat = "@"
$ws.Range("C61").Value = "1"
$ws.Range("C61").NumberFormat = "General"
$ws.Range("H61").Value = 1

# --- Row 62: LED Integration and function -----------------------------------
$ws.Range("A62").Value = 43207
$ws.Range("B62").Value = "LED Integration and function"
$ws.Range("C62").Value = 2
$ws.Range("F62").Value = 1

# --- Row 63: Led VHDL ---------------------------------------------------
# Hours (C63) is also a text "2" value in the source workbook.
$ws.Range("A63").Value = 43207
$ws.Range("B63").Value = "Led VHDL"
$ws.Range("C63").NumberFormat = "@"
$ws.Range("C63").Value = "2"
$ws.Range("C63").NumberFormat = "General"
$ws.Range("F63").Value = 1

# --- Row 64: Verify C function, wrapper, and functional test logic ---------
$ws.Range("A64").Value = 43207
$ws.Range("B64").Value = "Verify C function, wrapper, and functional test logic"
$ws.Range("C64").Value = 2
$ws.Range("I64").Value = 1

# Row heights settle (wrap-text autofit) after the new rows of task text are
# entered - match the resulting layout from the source edit.
$ws.Rows.Item(54).RowHeight = 26.25
$ws.Rows.Item(55).RowHeight = 15
$ws.Rows.Item(56).RowHeight = 15
$ws.Rows.Item(57).RowHeight = 26.25
$ws.Rows.Item(58).RowHeight = 26.25
$ws.Rows.Item(59).RowHeight = 26.25
$ws.Rows.Item(64).RowHeight = 13.8

# Column A widens slightly to accommodate the new dates.
$ws.Columns.Item(1).ColumnWidth = 10.5

# Leave the cursor where the author's session ended.
[void]$ws.Range("I65").Select()
